$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") width: 31 -> 13 characters.
# Excel's COM ColumnWidth setter stores (input + ~0.8333) chars in the XML,
# so back the input off by that padding to land on an exact 13.
$ws.Columns.Item(7).ColumnWidth = 12.1666667

# Replace every "Miss Dina Nasr" / "Miss Dina Nasr, Administrator" value in
# column G (rows 2..259) with the academic year "2025/2026".
$lastRow = 259
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*Dina Nasr*") {
        $cell.Value2 = "2025/2026"
    }
}
